$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (shared-string), never letting
# Excel auto-coerce a purely-numeric-looking string into a Number cell, and
# without leaving a NumberFormat/style behind on the destination cell.
#
# We stage the text (forced to Text format) in a scratch cell far outside the
# used range, copy it, and PasteSpecial only the *values* into the target -
# this carries the "this is text" flag along with the value but leaves the
# destination cell's style/format completely untouched.
# ---------------------------------------------------------------------------
$staging = $ws.Range("ZZ1")

function Set-TextValue($rangeAddress, $text) {
    $staging.NumberFormat = "@"
    $staging.Value = $text
    $staging.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
}

# Row 2: PickupID changes
Set-TextValue "C2" "10223102"

# Row 3: PickupID changes
Set-TextValue "C3" "10223105"

# Row 4: PickupID changes
Set-TextValue "C4" "10223388"

# Row 5: PickupID + Fail Log changes
Set-TextValue "C5" "10223226"
$ws.Range("F5").Value = "Cannot invoke ""org.openqa.selenium.WebElement.getText()"" because the return value of ""connect_OCBaseMethods.TCAcknowledge.isElementPresent(String)"" is null"

# The long Selenium "no such element" failure text shared by rows 11-14.
$noSuchElementMsg = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: headless chrome=119.0.6045.200)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.14', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '19.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 119.0.6045.200, chrome: {chromedriverVersion: 119.0.6045.105 (38c72552c5e..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:60496}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: b7e61e79fe2d4e962a4e6aeb1522a4a5
*** Element info: {Using=id, value=lblServiceID}
'@

# Row 11: PickupID + Fail Log changes
Set-TextValue "C11" "10223231"
$ws.Range("F11").Value = $noSuchElementMsg

# Row 12: PickupID + Fail Log changes
Set-TextValue "C12" "10223233"
$ws.Range("F12").Value = $noSuchElementMsg

# Row 13: PickupID + Fail Log changes
Set-TextValue "C13" "10223246"
$ws.Range("F13").Value = $noSuchElementMsg

# Row 14: PickupID + Fail Log changes
Set-TextValue "C14" "10223256"
$ws.Range("F14").Value = $noSuchElementMsg

# Row 24: PickupID change
Set-TextValue "C24" "136398862"

# Clean up the scratch cell so it leaves no trace in the saved workbook.
$staging.Clear()
